# Add, after the existing "Procuratore" run inside the checkbox-option
# table cell, a space run followed by an italic explanatory run:
#   " (allegare la procura, tranne nel caso in cui l'attribuzione
#   dell'incarico risulti dalla visura camerale)"
#
# Both new runs carry the same <w:rFonts w:cstheme="minorHAnsi"/> the
# surrounding text uses; the first keeps the paragraph's complex-script
# bold toggle (bCs) like "Procuratore" does, the second is italicised.

$d = $word.ActiveDocument

$range = $d.Content
$found = $range.Find.Execute("Procuratore", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $range.Collapse(0)

    # --- insert the separating space right after "Procuratore" ---
    $spaceStart = $range.End
    $range.InsertAfter(" ")
    $spaceEnd = $range.End

    # Briefly flag the new space run as italic so the engine doesn't
    # silently fold it back into the preceding "Procuratore" run (which
    # is not italicised); we immediately clear the flag again so the
    # space ends up matching the target plain/bCs-only formatting.
    $spaceRange = $d.Range($spaceStart, $spaceEnd)
    $spaceRange.Font.Italic = 1

    # --- insert the italicised explanatory text after the space ---
    $noteRange = $d.Range($spaceEnd, $spaceEnd)
    $noteRange.InsertAfter("(allegare la procura, tranne nel caso in cui l’attribuzione dell’incarico risulti dalla visura camerale)")
    $noteEnd = $noteRange.End
    $noteRange2 = $d.Range($spaceEnd, $noteEnd)
    $noteRange2.Font.Italic = 1

    # now clear the temporary italic flag on the space run
    $spaceRange2 = $d.Range($spaceStart, $spaceEnd)
    $spaceRange2.Font.Italic = 0

    Write-Host "Inserted procura note after 'Procuratore'."
} else {
    Write-Host "WARNING: 'Procuratore' not found; document left unchanged."
}
